# Append a new "Régule" movement row to the historique sheet, right after
# the last currently used row (row 25 -> new row 26).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.UsedRange.Rows.Count + 1

$ws.Cells.Item($newRow, 1).Value = "2025-06-02 08:37:32"
$ws.Cells.Item($newRow, 2).Value = "Verrouilleur 600mm 1E 255281"
$ws.Cells.Item($newRow, 3).Value = "Régule"
$ws.Cells.Item($newRow, 4).Value = 1
$ws.Cells.Item($newRow, 5).Value = 158
$ws.Cells.Item($newRow, 6).Value = 157

# Reference codes are stored as text (even though they look numeric), so
# force text entry with a leading apostrophe, then reset the cell style to
# "Normal" so no extra numeric-format / quote-prefix styling is left behind
# on the cell (matches how the other Reference cells in the sheet look).
$ws.Cells.Item($newRow, 7).Value = "'6755773992"
$ws.Cells.Item($newRow, 7).Style = "Normal"
